$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 0.96110082545191533
$ws.Range("D2").Value = 0.95417425971606762
$ws.Range("BP2").Value = 0.73958028198348091
$ws.Range("B3").Value = 0.9993806985796192
$ws.Range("C4").Value = 0.6897089149062996
$ws.Range("O4").Value = 0.6306890490012107
$ws.Range("X4").Value = 0.83285061571440422
$ws.Range("G5").Value = 0.90236873251853367
$ws.Range("D6").Value = 0.50987537930513738
$ws.Range("E6").Value = 0.96790450532032479
$ws.Range("P6").Value = 0.83626237631722988
$ws.Range("I7").Value = 0.89534281103700408
$ws.Range("G8").Value = 0.96961960219326793
$ws.Range("I8").Value = 0.84486099753548283
$ws.Range("R9").Value = 0.9689646703098671
$ws.Range("H10").Value = 0.7814579538484312
$ws.Range("K10").Value = 0.87837628519237487
$ws.Range("L10").Value = 0.72689287498372179
$ws.Range("I11").Value = 0.54453665191728451
$ws.Range("K12").Value = 0.95702964548208347
$ws.Range("N12").Value = 0.97657581575532681
$ws.Range("K13").Value = 0.95929240412832928
$ws.Range("O13").Value = 0.94316817674255304
$ws.Range("Q13").Value = 0.84194077122452671
$ws.Range("M14").Value = 0.76658833869072818
$ws.Range("O14").Value = 0.93349466844673978
$ws.Range("AA16").Value = 0.71247568239666892
$ws.Range("P17").Value = 0.92869094038679911
$ws.Range("S17").Value = 0.83576627319285901
$ws.Range("T18").Value = 0.88582829327409462
$ws.Range("R19").Value = 0.97883166225457008
$ws.Range("T19").Value = 0.89795791675590386
$ws.Range("U19").Value = 0.8054000139297941
$ws.Range("T21").Value = 0.90790536406553057
$ws.Range("BA21").Value = 0.76323522168976121
$ws.Range("U22").Value = 0.5715998246741758
$ws.Range("W22").Value = 0.94905808057291674
$ws.Range("X22").Value = 0.81226088759201343
$ws.Range("AQ22").Value = 0.61857053025015163
$ws.Range("BF22").Value = 0.60471422111549589
$ws.Range("J23").Value = 0.69018585111505093
$ws.Range("U23").Value = 0.75716139681334194
$ws.Range("X23").Value = 0.77465093798830087
$ws.Range("U24").Value = 0.9540349443120455
$ws.Range("X25").Value = 0.99645452154500691
$ws.Range("Z25").Value = 0.8079615595611358
$ws.Range("BK25").Value = 0.97705962673399571
$ws.Range("AB26").Value = 0.81883647395669446
$ws.Range("AE26").Value = 0.88720178553194806
$ws.Range("BP26").Value = 0.95775150455495583
$ws.Range("N27").Value = 0.97309020556738668
$ws.Range("Y27").Value = 0.86878149911686098
$ws.Range("Z27").Value = 0.99137243242122908
$ws.Range("AB27").Value = 0.83686444461217246
$ws.Range("BE27").Value = 0.94112180471556184
$ws.Range("BN28").Value = 0.86707015611994864
$ws.Range("AD29").Value = 0.98806466869329512
$ws.Range("AE29").Value = 0.68882013567223432
$ws.Range("AE30").Value = 0.86251675601035727
$ws.Range("AD32").Value = 0.94473135379287743
$ws.Range("AE32").Value = 0.86585775019135724
$ws.Range("AG32").Value = 0.93212283921747763
$ws.Range("AH33").Value = 0.8834072920279592
$ws.Range("AI33").Value = 0.81957903422662137
$ws.Range("AF34").Value = 0.83810126248477745
$ws.Range("AJ34").Value = 0.81546672257284625
$ws.Range("AW34").Value = 0.82463802952414811
$ws.Range("BJ34").Value = 0.95413161296338078
$ws.Range("BD35").Value = 0.91227600466458014
$ws.Range("AI36").Value = 0.91974262547946051
$ws.Range("AK36").Value = 0.97467421397113507
$ws.Range("AI37").Value = 0.65804560879105323
$ws.Range("AM37").Value = 0.58784562182625555
$ws.Range("AJ38").Value = 0.93261984623808325
$ws.Range("AK38").Value = 0.7852974139468849
$ws.Range("AM38").Value = 0.94149097448007457
$ws.Range("AN38").Value = 0.84297747445574434
$ws.Range("AM40").Value = 0.87364664287496185
$ws.Range("AO40").Value = 0.76990722605910134
$ws.Range("AM41").Value = 0.96878866109262463
$ws.Range("AQ41").Value = 0.91747332788033575
$ws.Range("AN42").Value = 0.76280899277762582
$ws.Range("O43").Value = 0.96610494848358663
$ws.Range("AP43").Value = 0.78148251995440643
$ws.Range("AS43").Value = 0.73837127340329323
$ws.Range("S44").Value = 0.62648608458314237
$ws.Range("BA44").Value = 0.81094651562313469
$ws.Range("AU45").Value = 0.7284274152256387
$ws.Range("AS46").Value = 0.71295661371690278
$ws.Range("AU46").Value = 0.59614326753838975
$ws.Range("AV46").Value = 0.71433639592132547
$ws.Range("AW47").Value = 0.78964005015078331
$ws.Range("AU48").Value = 0.93193480564700826
$ws.Range("AX48").Value = 0.85737619373868057
$ws.Range("K49").Value = 0.85436868324846849
$ws.Range("L49").Value = 0.90911445705451532
$ws.Range("AV49").Value = 0.71555289406732581
$ws.Range("AB51").Value = 0.87273342097111883
$ws.Range("AW51").Value = 0.61904950654733604
$ws.Range("AX51").Value = 0.65456070787996568
$ws.Range("AZ51").Value = 0.86173960284963924
$ws.Range("BA51").Value = 0.66924228682961251
$ws.Range("AX52").Value = 0.78344332852859366
$ws.Range("BA52").Value = 0.77528872204279176
$ws.Range("AP53").Value = 0.98712908211319239
$ws.Range("AZ54").Value = 0.72106438232388326
$ws.Range("BA54").Value = 0.78333493107258634
$ws.Range("BO54").Value = 0.81446974457504806
$ws.Range("BD55").Value = 0.89718752924341771
$ws.Range("BE55").Value = 0.60814732069605726
$ws.Range("BB56").Value = 0.85408760178493015
$ws.Range("BD57").Value = 0.81577339287256279
$ws.Range("BF57").Value = 0.69443975944224767
$ws.Range("AL58").Value = 0.94167424189157523
$ws.Range("BE59").Value = 0.99723983169028696
$ws.Range("BI59").Value = 0.68049116431687184
$ws.Range("T60").Value = 0.70934715116374081
$ws.Range("BG60").Value = 0.78017221071759457
$ws.Range("BJ60").Value = 0.75493597504164711
$ws.Range("C61").Value = 0.88299215152694344
$ws.Range("D61").Value = 0.84541825834657125
$ws.Range("BH61").Value = 0.77210418320379759
$ws.Range("BJ61").Value = 0.77074573892393705
$ws.Range("BK61").Value = 0.99219601337756191
$ws.Range("A62").Value = 0.81292159141504106
$ws.Range("BL63").Value = 0.83677396799202808
$ws.Range("BJ64").Value = 0.84057526078682521
$ws.Range("BM64").Value = 0.78791394032246032
$ws.Range("BN64").Value = 0.91290434563248946
$ws.Range("BN65").Value = 0.6364811824135852
$ws.Range("BO65").Value = 0.90527439212463912
$ws.Range("M66").Value = 0.69032237693478304
$ws.Range("BP66").Value = 0.91691981666982292
$ws.Range("AO67").Value = 0.97163872876892676
$ws.Range("BH67").Value = 0.76495050259725039
$ws.Range("BO68").Value = 0.81236872212774913
